$d = $word.ActiveDocument

# Each target paragraph currently holds its whole text in a single plain
# run; Word's background spell/grammar checker would normally split that
# run around the identifier-looking tokens and drop in w:proofErr markers
# (spellStart/spellEnd around underscore_names, gramStart/gramEnd around
# the dotted/paren tokens it treats as sentence-like). We reproduce that
# exact run/proofErr split by rewriting each paragraph's OOXML in place.

# --- "products(product_id, name, category_id)" ---
$targetText = 'products(product_id, name, category_id)'
$runPattern = '<w:r\b[^>]*><w:t\b[^>]*>products\(product_id, name, category_id\)</w:t></w:r>'
$newRuns = '<w:proofErr w:type="gramStart"/><w:r><w:t>products(</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/><w:r><w:t>product_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, name, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>category_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r>'
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $p = $d.Paragraphs.Item($i)
  $pText = $p.Range.Text.TrimEnd([char]13)
  if ($pText -eq $targetText) {
    $openXml = $p.Range.WordOpenXML
    if ($openXml -match $runPattern) {
      $openXml = $openXml -replace $runPattern, $newRuns
      $p.Range.InsertXML($openXml)
      $found = $true
    }
    break
  }
}
if (-not $found) {
  Write-Output "WARNING: paragraph not found/updated for: $targetText"
}

# --- "categories(category_id, category_name)" ---
$targetText = 'categories(category_id, category_name)'
$runPattern = '<w:r\b[^>]*><w:t\b[^>]*>categories\(category_id, category_name\)</w:t></w:r>'
$newRuns = '<w:proofErr w:type="gramStart"/><w:r><w:t>categories(</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/><w:r><w:t>category_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>category_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r>'
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $p = $d.Paragraphs.Item($i)
  $pText = $p.Range.Text.TrimEnd([char]13)
  if ($pText -eq $targetText) {
    $openXml = $p.Range.WordOpenXML
    if ($openXml -match $runPattern) {
      $openXml = $openXml -replace $runPattern, $newRuns
      $p.Range.InsertXML($openXml)
      $found = $true
    }
    break
  }
}
if (-not $found) {
  Write-Output "WARNING: paragraph not found/updated for: $targetText"
}

# --- "sales(sale_id, product_id, quantity, sale_date)" ---
$targetText = 'sales(sale_id, product_id, quantity, sale_date)'
$runPattern = '<w:r\b[^>]*><w:t\b[^>]*>sales\(sale_id, product_id, quantity, sale_date\)</w:t></w:r>'
$newRuns = '<w:proofErr w:type="gramStart"/><w:r><w:t>sales(</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/><w:r><w:t>sale_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>product_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, quantity, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sale_date</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r>'
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $p = $d.Paragraphs.Item($i)
  $pText = $p.Range.Text.TrimEnd([char]13)
  if ($pText -eq $targetText) {
    $openXml = $p.Range.WordOpenXML
    if ($openXml -match $runPattern) {
      $openXml = $openXml -replace $runPattern, $newRuns
      $p.Range.InsertXML($openXml)
      $found = $true
    }
    break
  }
}
if (-not $found) {
  Write-Output "WARNING: paragraph not found/updated for: $targetText"
}

# --- "    p.name AS product_name," ---
$targetText = '    p.name AS product_name,'
$runPattern = '<w:r\b[^>]*><w:t\b[^>]*>    p\.name AS product_name,</w:t></w:r>'
$newRuns = '<w:r><w:t xml:space="preserve">    p.name AS </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>product_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>,</w:t></w:r>'
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $p = $d.Paragraphs.Item($i)
  $pText = $p.Range.Text.TrimEnd([char]13)
  if ($pText -eq $targetText) {
    $openXml = $p.Range.WordOpenXML
    if ($openXml -match $runPattern) {
      $openXml = $openXml -replace $runPattern, $newRuns
      $p.Range.InsertXML($openXml)
      $found = $true
    }
    break
  }
}
if (-not $found) {
  Write-Output "WARNING: paragraph not found/updated for: $targetText"
}

# --- "    SUM(s.quantity) AS total_quantity_sold" ---
$targetText = '    SUM(s.quantity) AS total_quantity_sold'
$runPattern = '<w:r\b[^>]*><w:t\b[^>]*>    SUM\(s\.quantity\) AS total_quantity_sold</w:t></w:r>'
$newRuns = '<w:r><w:t xml:space="preserve">    SUM(</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>s.quantity</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">) AS </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>total_quantity_sold</w:t></w:r><w:proofErr w:type="spellEnd"/>'
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $p = $d.Paragraphs.Item($i)
  $pText = $p.Range.Text.TrimEnd([char]13)
  if ($pText -eq $targetText) {
    $openXml = $p.Range.WordOpenXML
    if ($openXml -match $runPattern) {
      $openXml = $openXml -replace $runPattern, $newRuns
      $p.Range.InsertXML($openXml)
      $found = $true
    }
    break
  }
}
if (-not $found) {
  Write-Output "WARNING: paragraph not found/updated for: $targetText"
}

# --- "    sales s ON p.product_id = s.product_id" ---
$targetText = '    sales s ON p.product_id = s.product_id'
$runPattern = '<w:r\b[^>]*><w:t\b[^>]*>    sales s ON p\.product_id = s\.product_id</w:t></w:r>'
$newRuns = '<w:r><w:t xml:space="preserve">    sales s ON </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>p.product</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>s.product</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_id</w:t></w:r><w:proofErr w:type="spellEnd"/>'
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $p = $d.Paragraphs.Item($i)
  $pText = $p.Range.Text.TrimEnd([char]13)
  if ($pText -eq $targetText) {
    $openXml = $p.Range.WordOpenXML
    if ($openXml -match $runPattern) {
      $openXml = $openXml -replace $runPattern, $newRuns
      $p.Range.InsertXML($openXml)
      $found = $true
    }
    break
  }
}
if (-not $found) {
  Write-Output "WARNING: paragraph not found/updated for: $targetText"
}

# --- "    p.product_id, p.name;" ---
$targetText = '    p.product_id, p.name;'
$runPattern = '<w:r\b[^>]*><w:t\b[^>]*>    p\.product_id, p\.name;</w:t></w:r>'
$newRuns = '<w:r><w:t xml:space="preserve">    </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>p.product</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, p.name;</w:t></w:r>'
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $p = $d.Paragraphs.Item($i)
  $pText = $p.Range.Text.TrimEnd([char]13)
  if ($pText -eq $targetText) {
    $openXml = $p.Range.WordOpenXML
    if ($openXml -match $runPattern) {
      $openXml = $openXml -replace $runPattern, $newRuns
      $p.Range.InsertXML($openXml)
      $found = $true
    }
    break
  }
}
if (-not $found) {
  Write-Output "WARNING: paragraph not found/updated for: $targetText"
}

